# A new weekly price record (row) for Jengibre at "Vega Modelo de Temuco"
# is inserted at row 68, pushing the existing rows 68-145 down to 69-146.
# The new row re-uses the same Mercado/Categoria/Unidad/Origen metadata as
# its neighbours, but carries its own Fecha/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68 (rows 68..145 shift down to 69..146)
$ws.Rows.Item(68).Insert()

# Seed the new row with the same static metadata columns as the row above
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion)
$ws.Range("A67:R67").Copy($ws.Range("A68:R68"))

# Now set the row-specific values for the newly inserted record
$ws.Range("D68").Value2 = 44638
$ws.Range("J68").Value2 = 50
$ws.Range("K68").Value2 = 20000
$ws.Range("L68").Value2 = 20000
$ws.Range("M68").Value2 = 20000
$ws.Range("P68").Value2 = 1538
